$d = $word.ActiveDocument

# Docx writer: Use different style for block quotes in notes.
# Add a new "Footnote Block Text" paragraph style, based on (and
# followed by) "Footnote Text", mirroring how "Block Text" is derived
# from "Body Text" -- so footnote block quotes can get their own font
# size independent of the ordinary body "Block Text" style.

$newStyle = $d.Styles.Add("Footnote Block Text", 1)

$newStyle.BaseStyle = $d.Styles("Footnote Text")
$newStyle.NextParagraphStyle = $d.Styles("Footnote Text")
$newStyle.Priority = 9
$newStyle.UnhideWhenUsed = $true
$newStyle.QuickStyle = $true

$newStyle.ParagraphFormat.SpaceBefore = 5
$newStyle.ParagraphFormat.SpaceAfter = 5
$newStyle.ParagraphFormat.FirstLineIndent = 0
$newStyle.ParagraphFormat.LeftIndent = 24
$newStyle.ParagraphFormat.RightIndent = 24
